# Fix dog vaccination rate parameter / restructure the "Initial parameter" table:
#  - drop the static R0_dog_to_dog row (it becomes a *calculated* value later on)
#  - correct the Human_population and Humans_per_free_roaming_dog inputs
#  - append a new "Calculated" block that derives Humans_per_km2,
#    Free_roaming_dog_population, Free_roaming_dogs_per_km2 and (finally) R0_dog_to_dog
#    from the corrected inputs, the last two rows pulling their label text from the
#    external "Define_program" model workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Correct two existing input values
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = 13125164      # Human_population
$ws.Range("C6").Value = 15.6          # Humans_per_free_roaming_dog

# ---------------------------------------------------------------------------
# 2. Remove the old static "R0_dog_to_dog" row (row 9) - it is replaced by a
#    calculated row near the bottom of the sheet.
# ---------------------------------------------------------------------------
$ws.Rows(9).Delete()

# ---------------------------------------------------------------------------
# 3. Append the new "Calculated" rows (now rows 14-17)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Calculated"
$ws.Range("B14").Value = "Humans_per_km2"
$ws.Range("C14").Formula = "=C3/C2"
$ws.Range("D14").Value = "Humans per km2"

$ws.Range("A15").Value = "Calculated"
$ws.Range("B15").Value = "Free_roaming_dog_population"
$ws.Range("C15").Formula = "=C3/C6"
$ws.Range("D15").Formula = "=[RabiesEcon.xlsm]Define_program!B14"

$ws.Range("A16").Value = "Calculated"
$ws.Range("B16").Value = "Free_roaming_dogs_per_km2"
$ws.Range("C16").Formula = "=C15/C2"
$ws.Range("D16").Formula = "=[RabiesEcon.xlsm]Define_program!B15"

$ws.Range("A17").Value = "Calculated"
$ws.Range("B17").Value = "R0_dog_to_dog"
$ws.Range("C17").Formula = "=0.34*LN(C16)"
$ws.Range("D17").Value = "Rabies R0 Dog to Dog"

# ---------------------------------------------------------------------------
# 4. Cosmetic view-state touch ups to mirror the authoring session
# ---------------------------------------------------------------------------
$ws.Range("D23").Select()

$ws.Columns(1).ColumnWidth = 28.44140625
$ws.Columns(2).ColumnWidth = 38.6640625
$ws.Columns(3).ColumnWidth = 20.5546875
$ws.Columns(4).ColumnWidth = 53.44140625
$ws.Columns(5).ColumnWidth = 20.5546875
